# Generate Report for Handoff
#
# A new handoff run produced a new source-file GUID and new xlf hash, so
# every cell that referenced the old guid/hash/timestamps needs to move to
# the new ones, on all three sheets (Overview, zh-cn, de-de), including the
# hyperlink "display" text that mirrors the file name shown in A2/B2.

$wb = $excel.ActiveWorkbook

$oldGuid = "77ef109b-63e3-49af-a4d0-fba0a7daa72b"
$newGuid = "b06a409c-b6c0-4eb4-9e1f-bcb3cefc4ab8"

$oldHash = "81fedf8989624220fb8d0a4a55d32944f4ec9a20"
$newHash = "1b8ca57ada2f210b73e2cf31355322f290acd9c2"

$oldMdName = "$oldGuid.md"
$newMdName = "$newGuid.md"
$oldMdPath = "e2e\$oldGuid.md"
$newMdPath = "e2e\$newGuid.md"

$oldZhXlf = "$oldGuid.$oldHash.zh-cn.xlf"
$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"
$oldDeXlf = "$oldGuid.$oldHash.de-de.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

$oldGenDate = "2016-08-24 13:01:51"
$newGenDate = "2016-08-24 13:02:18"

$oldHandoffDate = "2016-08-24 13:01:46"
$newHandoffDate = "2016-08-24 13:02:11"

# The external hyperlink target (stored in the worksheet's .rels file) is
# unchanged by this commit - only the human readable "display" text moves
# to the new file name.
$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/94cb04b34492c219b7760da7f92cd0813714a4a0/e2e/$oldGuid.md"

# Recreating a hyperlink (this runtime only supports add/replace, not true
# in-place edits of the display text) resets the cell to the generic
# built-in "Hyperlink" look, so restore the workbook's original custom
# hyperlink font (underline + cornflower blue, RGB 100,149,237) afterwards.
function Set-HyperlinkDisplay($ws, $cellAddr, $display) {
    $range = $ws.Range($cellAddr)
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($range, $hyperlinkTarget, "", "", $display)
    $range.Font.Underline = 2
    $range.Font.Color = 15570276
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newMdPath
$wsOverview.Range("G2").Value = $newGenDate
Set-HyperlinkDisplay $wsOverview "B2" $newMdPath

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newHandoffDate
Set-HyperlinkDisplay $wsZh "A2" $newMdName

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("G2").Value = $newDeXlf
# "Latest Handback DateTime" on this sheet shares the same underlying
# string as the Overview sheet's "Latest HO Xliff Generate Date" cell, so
# it moves to the same new timestamp value.
$wsDe.Range("H2").Value = $newGenDate
Set-HyperlinkDisplay $wsDe "A2" $newMdName
